$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 1.57
$ws.Range("H2").Value2 = 3.6
$ws.Range("I2").Value2 = 6.25
$ws.Range("J2").Value2 = 2.25
$ws.Range("K2").Value2 = 2.05
$ws.Range("L2").Value2 = 7
$ws.Range("O2").Value2 = 1.44
$ws.Range("P2").Value2 = 2.63
$ws.Range("Q2").Value2 = 2.35
$ws.Range("R2").Value2 = 1.57
$ws.Range("S2").Value2 = 1.5
$ws.Range("T2").Value2 = 2.5
$ws.Range("U2").Value2 = 2.38
$ws.Range("V2").Value2 = 1.53
$ws.Range("X2").Value2 = 6
$ws.Range("Y2").Value2 = 9
$ws.Range("Z2").Value2 = 11
$ws.Range("AC2").Value2 = 7
$ws.Range("AF2").Value2 = 81
$ws.Range("AG2").Value2 = 12
$ws.Range("AH2").Value2 = 29
$ws.Range("AJ2").Value2 = 81
$ws.Range("AN2").Value2 = 3.25
$ws.Range("AO2").Value2 = 8.5
$ws.Range("AQ2").Value2 = 29
$ws.Range("AR2").Value2 = 51
$ws.Range("AT2").Value2 = 2.5
$ws.Range("AW2").Value2 = 7.5
$ws.Range("AX2").Value2 = 41
$ws.Range("G3").Value2 = 1.67
$ws.Range("I3").Value2 = 5.25
$ws.Range("Q3").Value2 = 2.2
$ws.Range("R3").Value2 = 1.65
$ws.Range("S3").Value2 = 1.41
$ws.Range("T3").Value2 = 2.62
$ws.Range("Z3").Value2 = 12
$ws.Range("AC3").Value2 = 8.5
$ws.Range("AQ3").Value2 = 29
$ws.Range("AT3").Value2 = 2.63
$ws.Range("AZ3").Value2 = 126
$ws.Range("S5").Value2 = 1.62
$ws.Range("K8").Value2 = 2.15
$ws.Range("M8").Value2 = 9.800000000000001
$ws.Range("N8").Value2 = 1.02
$ws.Range("O8").Value2 = 1.24
$ws.Range("P8").Value2 = 3.3
$ws.Range("Q8").Value2 = 1.78
$ws.Range("R8").Value2 = 1.93
$ws.Range("V8").Value2 = 2
$ws.Range("Y8").Value2 = 12.5
$ws.Range("AC8").Value2 = 11
$ws.Range("AF8").Value2 = 55
$ws.Range("AG8").Value2 = 7.9
$ws.Range("AH8").Value2 = 9.5
$ws.Range("AJ8").Value2 = 16.5
$ws.Range("AK8").Value2 = 14
$ws.Range("AL8").Value2 = 23
$ws.Range("AM8").Value2 = 400
$ws.Range("AT8").Value2 = 2.65
$ws.Range("AU8").Value2 = 6.8
$ws.Range("AV8").Value2 = 55
$ws.Range("AW8").Value2 = 3.8
$ws.Range("AY8").Value2 = 16.5
$ws.Range("G10").Value2 = 3.4
$ws.Range("N10").Value2 = 9
$ws.Range("Q10").Value2 = 2.15
$ws.Range("R10").Value2 = 1.67
$ws.Range("AJ10").Value2 = 19
$ws.Range("AP10").Value2 = 29
$ws.Range("AX10").Value2 = 12
$ws.Range("O11").Value2 = 1.3
$ws.Range("P11").Value2 = 3.4
$ws.Range("Q11").Value2 = 1.98
$ws.Range("R11").Value2 = 1.83
$ws.Range("G14").Value2 = 1.83
$ws.Range("I14").Value2 = 4
$ws.Range("Q14").Value2 = 2.1
$ws.Range("R14").Value2 = 1.7
$ws.Range("X14").Value2 = 8.5
$ws.Range("Y14").Value2 = 9
$ws.Range("AA14").Value2 = 17
$ws.Range("AJ14").Value2 = 41
$ws.Range("AK14").Value2 = 34
$ws.Range("G20").Value2 = 1.85
$ws.Range("I20").Value2 = 4.75
$ws.Range("J20").Value2 = 2.6
$ws.Range("U20").Value2 = 2
$ws.Range("V20").Value2 = 1.73
$ws.Range("X20").Value2 = 8
$ws.Range("AE20").Value2 = 17
$ws.Range("AH20").Value2 = 23
$ws.Range("AI20").Value2 = 17
$ws.Range("AO20").Value2 = 10
$ws.Range("AQ20").Value2 = 34
$ws.Range("AR20").Value2 = 51
$ws.Range("BA20").Value2 = 126
$ws.Range("M22").Value2 = 1.07
$ws.Range("N22").Value2 = 9
$ws.Range("O22").Value2 = 1.33
$ws.Range("P22").Value2 = 3.25
$ws.Range("G23").Value2 = 2.8
$ws.Range("I23").Value2 = 2.3
$ws.Range("J23").Value2 = 3.2
$ws.Range("AA23").Value2 = 19
$ws.Range("AH23").Value2 = 15
$ws.Range("H24").Value2 = 3.75
$ws.Range("I24").Value2 = 4.05
$ws.Range("AD24").Value2 = 7.5
$ws.Range("AE24").Value2 = 14.5
$ws.Range("AF24").Value2 = 60
$ws.Range("AK24").Value2 = 37
$ws.Range("AP24").Value2 = 17.5
$ws.Range("AU24").Value2 = 7.4
$ws.Range("G26").Value2 = 2.65
$ws.Range("H26").Value2 = 2.7
$ws.Range("I26").Value2 = 2.92
$ws.Range("K26").Value2 = 1.85
$ws.Range("L26").Value2 = 3.55
$ws.Range("N26").Value2 = 5.8
$ws.Range("Q26").Value2 = 2.4
$ws.Range("S26").Value2 = 1.52
$ws.Range("T26").Value2 = 2.22
$ws.Range("U26").Value2 = 1.93
$ws.Range("W26").Value2 = 6.5
$ws.Range("X26").Value2 = 12
$ws.Range("Y26").Value2 = 10.25
$ws.Range("AA26").Value2 = 28
$ws.Range("AB26").Value2 = 45
$ws.Range("AC26").Value2 = 6.1
$ws.Range("AD26").Value2 = 5.4
$ws.Range("AE26").Value2 = 16
$ws.Range("AF26").Value2 = 100
$ws.Range("AG26").Value2 = 7
$ws.Range("AK26").Value2 = 30
$ws.Range("AN26").Value2 = 4.35
$ws.Range("AP26").Value2 = 25
$ws.Range("AT26").Value2 = 2.2
$ws.Range("AU26").Value2 = 7.1
$ws.Range("AV26").Value2 = 75
$ws.Range("AW26").Value2 = 4.65
$ws.Range("AX26").Value2 = 17
$ws.Range("BA26").Value2 = 120
$ws.Range("H27").Value2 = 2.8
$ws.Range("I27").Value2 = 3.9
$ws.Range("Y27").Value2 = 9
$ws.Range("Z27").Value2 = 20
$ws.Range("AD27").Value2 = 5.6
$ws.Range("AE27").Value2 = 16
$ws.Range("AI27").Value2 = 13
$ws.Range("AY27").Value2 = 29
